$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 7859
$ws1.Range("F5").Value = 7859
$ws1.Range("F6").Value = 101
$ws1.Range("F8").Value = 2129
$ws1.Range("F9").Value = 8630
$ws1.Range("F13").Value = 5767
$ws1.Range("F15").Value = 2754
$ws1.Range("F20").Value = 608
$ws1.Range("F21").Value = 61
$ws1.Range("F22").Value = 3902
$ws1.Range("F23").Value = 80
$ws1.Range("F24").Value = 63
$ws1.Range("F25").Value = 58
$ws1.Range("F26").Value = 16
$ws1.Range("F27").Value = 176
$ws1.Range("F28").Value = 24
$ws1.Range("F29").Value = 5445
$ws1.Range("F30").Value = 7
$ws1.Range("F31").Value = 70
$ws1.Range("F33").Value = 390
$ws1.Range("F35").Value = 391
$ws1.Range("F36").Value = 2226
$ws1.Range("F37").Value = 1518
$ws1.Range("F39").Value = 1120
$ws1.Range("F40").Value = 4450
$ws1.Range("F41").Value = 77
$ws1.Range("F43").Value = 40
$ws1.Range("F44").Value = 3554
$ws1.Range("F46").Value = 2341
$ws1.Range("F50").Value = 19

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 122
$ws2.Range("F3").Value = 150
$ws2.Range("F5").Value = 70
$ws2.Range("F6").Value = 20

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 7859
$ws4.Range("F5").Value = 7859
$ws4.Range("F6").Value = 101
$ws4.Range("F8").Value = 2129
$ws4.Range("F9").Value = 8630
$ws4.Range("F13").Value = 5767
$ws4.Range("F15").Value = 2754
$ws4.Range("F19").Value = 122
$ws4.Range("F21").Value = 150
$ws4.Range("F22").Value = 608
$ws4.Range("F24").Value = 61
$ws4.Range("F25").Value = 3902
$ws4.Range("F26").Value = 80
$ws4.Range("F27").Value = 63
$ws4.Range("F28").Value = 58
$ws4.Range("F29").Value = 24
$ws4.Range("F30").Value = 5445
$ws4.Range("F31").Value = 70
$ws4.Range("F32").Value = 390
$ws4.Range("F34").Value = 391
$ws4.Range("F35").Value = 70
$ws4.Range("F36").Value = 2226
$ws4.Range("F37").Value = 1518
$ws4.Range("F38").Value = 20
$ws4.Range("F40").Value = 1120
$ws4.Range("F42").Value = 4450
$ws4.Range("F43").Value = 77
$ws4.Range("F45").Value = 40
$ws4.Range("F46").Value = 3554
$ws4.Range("F48").Value = 2341

Write-Output "done"
